# Add a new weekly price record for "Femacal de La Calera" (Haba).
# The new observation is inserted as row 8, pushing the existing rows
# 8..78 down to 9..79 (matching the target OOXML diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 8 (shifts rows 8-78 down to 9-79).
$ws.Rows.Item(8).Insert()

# Populate the new row 8 with the new weekly price observation.
$ws.Range("A8").Value = 3
$ws.Range("B8").Value = "Femacal de La Calera"
$ws.Range("C8").Value = "Coquimbo"
$ws.Range("D8").Value = 44490
$ws.Range("E8").Value = 5
$ws.Range("F8").Value = 100112026
$ws.Range("G8").Value = "Haba"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 110
$ws.Range("K8").Value = 9000
$ws.Range("L8").Value = 9500
$ws.Range("M8").Value = 9273
$ws.Range("N8").Value = "$/malla 25 kilos"
$ws.Range("O8").Value = "Provincia de Quillota"
$ws.Range("P8").Value = 371
$ws.Range("Q8").Value = 25
$ws.Range("R8").Value = "Hortaliza"
